$d = $word.ActiveDocument

# The cover-sheet table is the first table in the document. Append a new
# row at the end ("Report word count" / "3,269"), mirroring the existing
# label/value row pattern (bold label in column 1, plain value in column 2).
$t = $d.Tables(1)
$newRow = $t.Rows.Add()
$newRow.Cells(1).Range.Text = "Report word count"
$newRow.Cells(2).Range.Text = "3,269"
